# part of storeItemDefine rebuild
# Rebuild Sheet1 of SpecialItemDefine: drop the "level", "iconResource" and
# "desc" columns, keeping only ID + title (old column C becomes the new
# column B), and change the month_card row to reference the plain
# "month_card" string instead of the numeric level/extra icon/desc columns.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "level" column (B) entirely - remaining columns shift left.
$ws.Range("B1:B2").EntireColumn.Delete()

# Remove the "iconResource" column (now C after the previous delete) -
# remaining columns shift left again.
$ws.Range("C1:C2").EntireColumn.Delete()

# Remove the "desc" column (now C again) - leaves just ID / title.
$ws.Range("C1:C2").EntireColumn.Delete()

# Old "title" column (now B) keeps its header; fix up the data row so the
# month card item just references the shared "month_card" string.
$ws.Range("B1").Value = "title"
$ws.Range("B2").Value = "month_card"

$ws.Range("D1").Select()
